$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values in row 1 for new columns P and Q (continuing the 0..15 sequence),
# copying O1's format (border/bold/alignment) onto the new cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For data rows 2 through 25: update columns I, K, M, O and add P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
